$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Annotations")

# --- Append 7 new rows to the bottom of the Annotations table (rows 127-133) ---

# Row 127: Golden Bowl, The (Literary Work)
$ws.Range("A127").Value = "Golden Bowl, The"
$ws.Range("B127").Value = "Literary Work"
$ws.Range("C127").Value = "lit-gold"
$ws.Range("D127").Value = "../resources/annotations.xml#lit-gold"

# Row 128: Ambassadors, The (Literary Work)
$ws.Range("A128").Value = "Ambassadors, The"
$ws.Range("B128").Value = "Literary Work"
$ws.Range("C128").Value = "lit-amba"
$ws.Range("D128").Value = "../resources/annotations.xml#lit-amba"

# Row 129: Edward IV (Person) with note
$ws.Range("A129").Value = "Edward IV"
$ws.Range("B129").Value = "Person"
$ws.Range("C129").Value = "psn-edw4"
$ws.Range("D129").Value = "../resources/annotations.xml#psn-edw4"
$ws.Range("E129").Value = "Is this actually Edward IV??? I'm guessing based on Latimer and the c16th but then again… Check against quotes! "

# Row 130: Shepherd's Calendar, The (Literary Work)
$ws.Range("A130").Value = "Shepherd's Calendar, The"
$ws.Range("B130").Value = "Literary Work"
$ws.Range("C130").Value = "lit-shep"
$ws.Range("D130").Value = "../resources/annotations.xml#lit-shep"

# Row 131: Amesbury (Place)
$ws.Range("A131").Value = "Amesbury"
$ws.Range("B131").Value = "Place"
$ws.Range("C131").Value = "pla-ames"
$ws.Range("D131").Value = "../resources/annotations.xml#pla-ames"

# Row 132: Stonehenge (Place)
$ws.Range("A132").Value = "Stonehenge"
$ws.Range("B132").Value = "Place"
$ws.Range("C132").Value = "pla-ston"
$ws.Range("D132").Value = "../resources/annotations.xml#pla-ston"

# Row 133: Lundy, Isle of (Place)
$ws.Range("A133").Value = "Lundy, Isle of"
$ws.Range("B133").Value = "Place"
$ws.Range("C133").Value = "pla-lund"
$ws.Range("D133").Value = "../resources/annotations.xml#pla-lund"

# Match the row height (32pt, same as the rest of the table) used by these new rows
$ws.Range("A127:E133").EntireRow.RowHeight = 32

# --- Move the active selection on the Annotations sheet to E125 ---
$ws.Range("E125").Select()
